$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Stand" (as-of) timestamp string in A4
$ws.Range("A4").Value = "Stand: 13.12.2022, 16:14"

# Fill in forecast consumption values (column D, "Gesamt (Netzlast) [MWh]")
# for rows 1547-1642, which previously held a placeholder "-" text value
# and now hold the actual numeric forecast, formatted like the surrounding cells.
$forecastValues = @{
    1547 = 13508.5
    1548 = 13315.0
    1549 = 13258.75
    1550 = 13222.75
    1551 = 12920.0
    1552 = 12860.25
    1553 = 12790.25
    1554 = 12725.75
    1555 = 12744.0
    1556 = 12733.0
    1557 = 12801.5
    1558 = 12801.75
    1559 = 12872.5
    1560 = 12884.5
    1561 = 12893.75
    1562 = 12921.0
    1563 = 13110.0
    1564 = 13212.0
    1565 = 13284.5
    1566 = 13383.0
    1567 = 13760.75
    1568 = 13842.5
    1569 = 14113.25
    1570 = 14380.75
    1571 = 15087.75
    1572 = 15228.75
    1573 = 15528.0
    1574 = 15814.5
    1575 = 16233.75
    1576 = 16291.0
    1577 = 16541.5
    1578 = 16822.75
    1579 = 17113.25
    1580 = 17155.75
    1581 = 17250.0
    1582 = 17251.0
    1583 = 17322.0
    1584 = 17218.25
    1585 = 17335.5
    1586 = 17446.75
    1587 = 17571.75
    1588 = 17543.5
    1589 = 17668.5
    1590 = 17768.0
    1591 = 17799.25
    1592 = 17825.75
    1593 = 17935.25
    1594 = 17936.5
    1595 = 17875.0
    1596 = 17929.25
    1597 = 17882.0
    1598 = 17795.5
    1599 = 17805.25
    1600 = 17630.0
    1601 = 17429.25
    1602 = 17330.25
    1603 = 17457.25
    1604 = 17361.0
    1605 = 17193.25
    1606 = 17113.25
    1607 = 17206.5
    1608 = 17157.25
    1609 = 17134.25
    1610 = 17073.75
    1611 = 17211.0
    1612 = 17223.25
    1613 = 17304.25
    1614 = 17449.25
    1615 = 17589.25
    1616 = 17701.75
    1617 = 17764.0
    1618 = 17738.0
    1619 = 17644.5
    1620 = 17581.0
    1621 = 17549.25
    1622 = 17544.75
    1623 = 17350.25
    1624 = 17189.0
    1625 = 16986.5
    1626 = 16841.25
    1627 = 16580.5
    1628 = 16234.25
    1629 = 15951.0
    1630 = 15829.5
    1631 = 15604.5
    1632 = 15544.0
    1633 = 15376.25
    1634 = 15217.75
    1635 = 15100.5
    1636 = 14986.0
    1637 = 14723.75
    1638 = 14584.0
    1639 = 14400.25
    1640 = 14269.5
    1641 = 14152.0
    1642 = 13945.0
}

foreach ($row in $forecastValues.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $forecastValues[$row]
    $cell.NumberFormat = "#,##0.00"
}
